$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Paper Language" header in column G, matching the style of
# the existing "Paper Name" header in column F.
$ws.Range("G1").Value = "Paper Language"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# New column G width, matching the authored width from the diff.
$ws.Columns.Item(7).ColumnWidth = 17.7

# Row 1 height shrinks now that wrapping is no longer forcing extra height.
$ws.Rows.Item(1).RowHeight = 15.75

# Update selection to match the authored state.
$ws.Range("G15").Select()
